# Khawase_LabExam03Grading.xlsx - "from 33-41 - Driver" grading update
#
# Fills in the grading scores for the "CustomerMappingDriver Class" section
# (rows 27-31) and the "Compilation errors if any" row (row 37):
#   - Row 29 (For successfully scanning data from input file): Total Points = 16
#   - Row 30 (For correct and properly aligned output): Total Points = 0,
#       with grading comment "For incorrect output"
#   - Row 37 (Compilation errors if any): Total Points = -5,
#       with grading comment "Compilation error"
# The Total rows (31, 38) contain SUM formulas and recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 29: scanning data from input file scored full 16 points ---
$ws.Range("E29").Value = 16

# --- Row 37: compilation errors deduction applied (-5), with comment ---
# (Written before row 30's comment so the shared-string table gets
# "Compilation error" before "For incorrect output", matching the order
# in which these grading notes were entered.)
$ws.Range("F37").Value = "Compilation error"
$ws.Range("E37").Value = -5

# --- Row 30: output scored 0 points, with explanatory comment ---
$ws.Range("F30").Value = "For incorrect output"
$ws.Range("E30").Value = 0

# --- Update the view/selection to reflect where the grader was working ---
$ws.Activate()
$ws.Range("F30").Select()
$excel.ActiveWindow.ScrollRow = 21
$excel.ActiveWindow.ScrollColumn = 1
